$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.74209299016793
$ws.Range("C2").Value = 8.620208856397925
$ws.Range("D2").Value = 7.148593320140064
$ws.Range("F2").Value = 42.85088131713091
$ws.Range("G2").Value = 51.79652735140372
$ws.Range("H2").Value = 19.94904172104233
$ws.Range("J2").Value = 10.64650781202987
$ws.Range("L2").Value = 11.4838937008914
$ws.Range("M2").Value = 17.99075101424994
$ws.Range("B3").Value = 19.40185026025539
$ws.Range("C3").Value = 8.317108092986068
$ws.Range("D3").Value = 7.148507805370994
$ws.Range("F3").Value = 42.90088622264653
$ws.Range("G3").Value = 51.77207886595306
$ws.Range("H3").Value = 19.99670212024472
$ws.Range("J3").Value = 10.66547099276632
$ws.Range("L3").Value = 11.49298133412175
$ws.Range("M3").Value = 17.93155874624999
$ws.Range("B4").Value = 19.19419295380378
$ws.Range("C4").Value = 8.12367075531915
$ws.Range("D4").Value = 7.14865803923872
$ws.Range("F4").Value = 42.94224469423641
$ws.Range("G4").Value = 51.77298730876744
$ws.Range("H4").Value = 20.02997847157386
$ws.Range("J4").Value = 10.67769557892754
$ws.Range("L4").Value = 11.49985369900504
$ws.Range("M4").Value = 17.89806912747
$ws.Range("B5").Value = 19.10999952614327
$ws.Range("C5").Value = 8.043070345135087
$ws.Range("D5").Value = 7.148770397084945
$ws.Range("F5").Value = 42.96177010066383
$ws.Range("G5").Value = 51.77735266598651
$ws.Range("H5").Value = 20.04454527452134
$ws.Range("J5").Value = 10.68282379717229
$ws.Range("L5").Value = 11.50297986467679
$ws.Range("M5").Value = 17.88514917167204
$ws.Range("B6").Value = 19.09604848768011
$ws.Range("C6").Value = 8.029581839350911
$ws.Range("D6").Value = 7.148792147130316
$ws.Range("F6").Value = 42.96517337962008
$ws.Range("G6").Value = 51.77831852626186
$ws.Range("H6").Value = 20.0470247876321
$ws.Range("J6").Value = 10.68368420211261
$ws.Range("L6").Value = 11.50351864577504
$ws.Range("M6").Value = 17.88304801635478
$ws.Range("B7").Value = 19.1930556059603
$ws.Range("C7").Value = 8.122590830020911
$ws.Range("D7").Value = 7.148659347264948
$ws.Range("F7").Value = 42.94249721524798
$ws.Range("G7").Value = 51.77303001855855
$ws.Range("H7").Value = 20.03017085321512
$ws.Range("J7").Value = 10.67776414566861
$ws.Range("L7").Value = 11.49989454036932
$ws.Range("M7").Value = 17.89789192745029
$ws.Range("B8").Value = 19.62458504769688
$ws.Range("C8").Value = 8.51726141753883
$ws.Range("D8").Value = 7.148521847843289
$ws.Range("F8").Value = 42.86590710996005
$ws.Range("G8").Value = 51.78478952467876
$ws.Range("H8").Value = 19.96464067200375
$ws.Range("J8").Value = 10.652925992088
$ws.Range("L8").Value = 11.48675922765473
$ws.Range("M8").Value = 17.96975493573976
$ws.Range("B9").Value = 20.47582937025805
$ws.Range("C9").Value = 9.230327256638862
$ws.Range("D9").Value = 7.149854545404202
$ws.Range("F9").Value = 42.8005642595525
$ws.Range("G9").Value = 51.93437313925764
$ws.Range("H9").Value = 19.86808657720339
$ws.Range("J9").Value = 10.60880756720757
$ws.Range("L9").Value = 11.4712292683082
$ws.Range("M9").Value = 18.13288265427936
$ws.Range("B10").Value = 21.09785104025691
$ws.Range("C10").Value = 9.713995437287458
$ws.Range("D10").Value = 7.151802777890397
$ws.Range("F10").Value = 42.80463292701344
$ws.Range("G10").Value = 52.12146233953975
$ws.Range("H10").Value = 19.81677179566998
$ws.Range("J10").Value = 10.57916068450203
$ws.Range("L10").Value = 11.46601738721548
$ws.Range("M10").Value = 18.26563202330193
$ws.Range("B11").Value = 21.37878678416142
$ws.Range("C11").Value = 9.924770278064038
$ws.Range("D11").Value = 7.152897942353488
$ws.Range("F11").Value = 42.81784075340988
$ws.Range("G11").Value = 52.22327158595325
$ws.Range("H11").Value = 19.79771654046116
$ws.Range("J11").Value = 10.56626779927419
$ws.Range("L11").Value = 11.46498340928485
$ws.Range("M11").Value = 18.32867581205278
$ws.Range("B12").Value = 21.4847756054956
$ws.Range("C12").Value = 10.00321907170739
$ws.Range("D12").Value = 7.153342551482607
$ws.Range("F12").Value = 42.82447702683634
$ws.Range("G12").Value = 52.26421415702539
$ws.Range("H12").Value = 19.79111944375024
$ws.Range("J12").Value = 10.5614704676503
$ws.Range("L12").Value = 11.46478327934303
$ws.Range("M12").Value = 18.35291724963547
$ws.Range("B13").Value = 21.46196821211896
$ws.Range("C13").Value = 9.986385058654285
$ws.Range("D13").Value = 7.153245469715451
$ws.Range("F13").Value = 42.82297506888697
$ws.Range("G13").Value = 52.25529039382069
$ws.Range("H13").Value = 19.79251269770748
$ws.Range("J13").Value = 10.56249988878636
$ws.Range("L13").Value = 11.46481788147009
$ws.Range("M13").Value = 18.34768027384494
$ws.Range("B14").Value = 21.38751505348062
$ws.Range("C14").Value = 9.931251900291919
$ws.Range("D14").Value = 7.15293392232164
$ws.Range("F14").Value = 42.81835395616903
$ws.Range("G14").Value = 52.22659214649516
$ws.Range("H14").Value = 19.79716138205801
$ws.Range("J14").Value = 10.56587142041245
$ws.Range("L14").Value = 11.46496311318157
$ws.Range("M14").Value = 18.33066287239267
$ws.Range("B15").Value = 21.34185583382533
$ws.Range("C15").Value = 9.897302201058613
$ws.Range("D15").Value = 7.152746979338664
$ws.Range("F15").Value = 42.81573630991254
$ws.Range("G15").Value = 52.20932441987942
$ws.Range("H15").Value = 19.80008946930394
$ws.Range("J15").Value = 10.56794762753584
$ws.Range("L15").Value = 11.46507697410413
$ws.Range("M15").Value = 18.32028674710612
$ws.Range("B16").Value = 21.07944145873858
$ws.Range("C16").Value = 9.70003141060109
$ws.Range("D16").Value = 7.151735397481677
$ws.Range("F16").Value = 42.80399846699864
$ws.Range("G16").Value = 52.11514399837328
$ws.Range("H16").Value = 19.81810357189599
$ws.Range("J16").Value = 10.58001517263276
$ws.Range("L16").Value = 11.46611179240224
$ws.Range("M16").Value = 18.2615642360575
$ws.Range("B17").Value = 20.91786439808208
$ws.Range("C17").Value = 9.576615352453215
$ws.Range("D17").Value = 7.151168225700419
$ws.Range("F17").Value = 42.79970807162454
$ws.Range("G17").Value = 52.06163762510245
$ws.Range("H17").Value = 19.83025447123745
$ws.Range("J17").Value = 10.58756995032406
$ws.Range("L17").Value = 11.4670884622377
$ws.Range("M17").Value = 18.22621039776836
$ws.Range("B18").Value = 20.82474488660629
$ws.Range("C18").Value = 9.504761084536094
$ws.Range("D18").Value = 7.15086167513814
$ws.Range("F18").Value = 42.79830931000071
$ws.Range("G18").Value = 52.03243515062478
$ws.Range("H18").Value = 19.83764684957989
$ws.Range("J18").Value = 10.59197116003357
$ws.Range("L18").Value = 11.46777610833374
$ws.Range("M18").Value = 18.20612660272625
$ws.Range("B19").Value = 20.79318775919423
$ws.Range("C19").Value = 9.480284497049691
$ws.Range("D19").Value = 7.150761265486609
$ws.Range("F19").Value = 42.79801923588789
$ws.Range("G19").Value = 52.02281815988535
$ws.Range("H19").Value = 19.84021901698967
$ws.Range("J19").Value = 10.59347094855543
$ws.Range("L19").Value = 11.46803058007971
$ws.Range("M19").Value = 18.19937005924767
$ws.Range("B20").Value = 20.93508443412164
$ws.Range("C20").Value = 9.589843386578394
$ws.Range("D20").Value = 7.151226567206814
$ws.Range("F20").Value = 42.80005414228011
$ws.Range("G20").Value = 52.06717074847135
$ws.Range("H20").Value = 19.8289192106934
$ws.Range("J20").Value = 10.58675994896408
$ws.Range("L20").Value = 11.46697147068662
$ws.Range("M20").Value = 18.22994801545065
$ws.Range("B21").Value = 21.40939528997735
$ws.Range("C21").Value = 9.947483217995785
$ws.Range("D21").Value = 7.153024621182507
$ws.Range("F21").Value = 42.8196669159184
$ws.Range("G21").Value = 52.23495677807961
$ws.Range("H21").Value = 19.79577914328021
$ws.Range("J21").Value = 10.56487881851417
$ws.Range("L21").Value = 11.46491526704637
$ws.Range("M21").Value = 18.3356514168423
$ws.Range("B22").Value = 21.71703503120389
$ws.Range("C22").Value = 10.17323948815719
$ws.Range("D22").Value = 7.154373932506375
$ws.Range("F22").Value = 42.84201339578777
$ws.Range("G22").Value = 52.35853498643288
$ws.Range("H22").Value = 19.7777273754086
$ws.Range("J22").Value = 10.55107303626993
$ws.Range("L22").Value = 11.46468671130753
$ws.Range("M22").Value = 18.406873689617
$ws.Range("B23").Value = 21.55309008785487
$ws.Range("C23").Value = 10.0534904572127
$ws.Range("D23").Value = 7.153637890077272
$ws.Range("F23").Value = 42.8292146240089
$ws.Range("G23").Value = 52.29131018972774
$ws.Range("H23").Value = 19.78703124351928
$ws.Range("J23").Value = 10.5583963120965
$ws.Range("L23").Value = 11.46470693462986
$ws.Range("M23").Value = 18.36866996721702
$ws.Range("B24").Value = 20.92729995398018
$ws.Range("C24").Value = 9.583865795102902
$ws.Range("D24").Value = 7.151200130198124
$ws.Range("F24").Value = 42.79989435763984
$ws.Range("G24").Value = 52.06466436624668
$ws.Range("H24").Value = 19.82952161540683
$ws.Range("J24").Value = 10.58712597029461
$ws.Range("L24").Value = 11.46702396958175
$ws.Range("M24").Value = 18.22825748491295
$ws.Range("B25").Value = 20.24568501120472
$ws.Range("C25").Value = 9.044288962812946
$ws.Range("D25").Value = 7.149323278067935
$ws.Range("F25").Value = 42.80911645025295
$ws.Range("G25").Value = 51.88034113479129
$ws.Range("H25").Value = 19.89077083960289
$ws.Range("J25").Value = 10.62025465975893
$ws.Range("L25").Value = 11.47433921130768
$ws.Range("M25").Value = 18.08643837747508
